# Update "想去人数" (F column) figures across sheets as generated at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 198
$ws1.Range("F4").Value  = 60
$ws1.Range("F5").Value  = 1638
$ws1.Range("F6").Value  = 3244
$ws1.Range("F7").Value  = 791
$ws1.Range("F8").Value  = 2005
$ws1.Range("F9").Value  = 1929
$ws1.Range("F10").Value = 987
$ws1.Range("F11").Value = 345
$ws1.Range("F13").Value = 1593
$ws1.Range("F14").Value = 341
$ws1.Range("F17").Value = 34
$ws1.Range("F18").Value = 1408
$ws1.Range("F19").Value = 503
$ws1.Range("F20").Value = 606
$ws1.Range("F21").Value = 307
$ws1.Range("F22").Value = 10559
$ws1.Range("F23").Value = 9726
$ws1.Range("F25").Value = 647
$ws1.Range("F26").Value = 1814
$ws1.Range("F27").Value = 136
$ws1.Range("F28").Value = 395

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 115

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 198
$ws4.Range("F6").Value  = 60
$ws4.Range("F7").Value  = 1638
$ws4.Range("F8").Value  = 3244
$ws4.Range("F9").Value  = 791
$ws4.Range("F10").Value = 2005
$ws4.Range("F11").Value = 1929
$ws4.Range("F12").Value = 987
$ws4.Range("F13").Value = 345
$ws4.Range("F15").Value = 1593
$ws4.Range("F16").Value = 341
$ws4.Range("F21").Value = 34
$ws4.Range("F22").Value = 1408
$ws4.Range("F23").Value = 503
$ws4.Range("F24").Value = 606
$ws4.Range("F25").Value = 307
$ws4.Range("F26").Value = 10559
$ws4.Range("F27").Value = 9726
$ws4.Range("F29").Value = 647
$ws4.Range("F30").Value = 1814
$ws4.Range("F31").Value = 115
$ws4.Range("F33").Value = 136
$ws4.Range("F34").Value = 395
